$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before D, shifting existing D:K data to E:L.
$ws.Columns("D").Insert()

# Copy number formats/styles from the (now shifted) column E into the new column D
# so the new column matches its siblings (date format row 7/38/80, number format elsewhere).
$srcRange = $ws.Range("E5:E102")
$dstRange = $ws.Range("D5:D102")
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Match the new column's width to its neighbor.
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# Populate the new column D with the latest period's figures.
$ws.Range("D7").Value = 43312
$ws.Range("D8").Value = 645300
$ws.Range("D9").Value = 544000
$ws.Range("D10").Value = 101300
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 300
$ws.Range("D15").Value = 20300
$ws.Range("D17").Value = 653600
$ws.Range("D18").Value = -8300
$ws.Range("D20").Value = 2900
$ws.Range("D21").Value = 31700
$ws.Range("D22").Value = 29900
$ws.Range("D23").Value = -35300
$ws.Range("D24").Value = -71200
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 35900
$ws.Range("D27").Value = 35400
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -2900
$ws.Range("D33").Value = 35400
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 35400
$ws.Range("D38").Value = 43312
$ws.Range("D41").Value = 92100
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 99300
$ws.Range("D44").Value = 47800
$ws.Range("D45").Value = 25100
$ws.Range("D46").Value = 264300
$ws.Range("D47").Value = "NA"
$ws.Range("D48").Value = 106600
$ws.Range("D49").Value = 447300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 8800
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 827000
$ws.Range("D57").Value = 78200
$ws.Range("D58").Value = 70300
$ws.Range("D59").Value = 142100
$ws.Range("D60").Value = 290600
$ws.Range("D61").Value = 383100
$ws.Range("D62").Value = 10500
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 684200
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 35200
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -7363600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 107600
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43312
$ws.Range("D81").Value = 35400
$ws.Range("D83").Value = 37100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 11800
$ws.Range("D91").Value = -18400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -452300
$ws.Range("D96").Value = -1100
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 421900
$ws.Range("D101").Value = 100
$ws.Range("D102").Value = -18500

Write-Output "Edit complete"
